# Mudanças p/ nova correção de atividade incompleta
#
# Slide 9 (old) was a "blank" slide with just a grey background and a
# picture (a placeholder for a future "Resultados" slide). This edit:
#   1. Turns the old slide 9 into a proper "Resultados" title/subtitle slide.
#   2. Inserts two new title/subtitle slides ("UX, UI e Design Interativo"
#      and "Dificuldades") right after it.
#   3. Moves the original picture content to a brand-new slide at the very
#      end of the deck.
#
# Final slide order: 1..8 (unchanged), 9=Resultados (new text),
# 10=UX, UI e Design Interativo, 11=Dificuldades, 12=old slide9 content.

$p = $ppt.ActivePresentation

$sectionLayout = $p.SlideMaster.CustomLayouts.Item(8)   # SECTION_TITLE_AND_DESCRIPTION

function Set-TitleSubtitle($slide, $titleText, $subtitleParts) {
    # Apply the SECTION_TITLE_AND_DESCRIPTION layout, then strip away the
    # "body" placeholder (idx=2) and any non-placeholder legacy shapes so
    # only title + subtitle remain, matching the target structure.
    $slide.FollowMasterBackground = -1
    $slide.CustomLayout = $sectionLayout

    for ($j = $slide.Shapes.Count; $j -ge 1; $j--) {
        $shp = $slide.Shapes.Item($j)
        $phType = 0
        try { $phType = $shp.PlaceholderFormat.Type } catch {}
        if ($phType -eq 2) {
            $shp.Delete()
        } elseif ($phType -ne 1 -and $phType -ne 4) {
            $shp.Delete()
        }
    }

    $title = $slide.Shapes.Item(1)
    $title.TextFrame.TextRange.Text = $titleText
    $title.TextFrame.TextRange.LanguageID = "pt-BR"

    $sub = $slide.Shapes.Item(2)
    $sub.TextFrame.TextRange.Text = $subtitleParts[0]
    $sub.TextFrame.TextRange.LanguageID = "pt-BR"
    for ($k = 1; $k -lt $subtitleParts.Count; $k++) {
        $appended = $sub.TextFrame.TextRange.InsertAfter($subtitleParts[$k])
        $appended.LanguageID = "pt-BR"
    }
}

# --- Create the two new slides first (so they claim slide10.xml/slide11.xml),
#     then duplicate the old slide 9 last (so it becomes slide12.xml). ---
$uxSlide = $p.Slides.AddSlide($p.Slides.Count + 1, $sectionLayout)
$diffSlide = $p.Slides.AddSlide($p.Slides.Count + 1, $sectionLayout)

$oldSlide9 = $p.Slides.Item(9)
$pictureSlide = $oldSlide9.Duplicate()
$pictureSlide.MoveTo($p.Slides.Count)

# --- Fill in the content for each slide ---
Set-TitleSubtitle $oldSlide9 "Resultados" @("Revisão das necessidades do briefing")
Set-TitleSubtitle $uxSlide "UX, UI e Design Interativo" @("De que maneira respondem os sites desenvolvidos a estes princípios")
Set-TitleSubtitle $diffSlide "Dificuldades" @("Dificuldades de ", "desenvolvimento para desktop, tablet e smartphone")
